$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "51.933.63"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3
$ws.Range("D3").Value = "2.915.98"
$ws.Range("E3").Value = "  +3.46%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
Set-TextValue $ws.Range("D5") "352.17"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
Set-TextValue $ws.Range("D6") "113.01"
$ws.Range("E6").Value = "  -0.19%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.558"
$ws.Range("E7").Value = "  -0.54%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.618"
$ws.Range("E9").Value = "  -0.46%  "

# Row 10
Set-TextValue $ws.Range("D10") "39.36"
$ws.Range("E10").Value = "  -2.62%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0872"
$ws.Range("E11").Value = "  +3.33%  "

# Row 12
$ws.Range("E12").Value = "  +0.64%  "

# Row 13
Set-TextValue $ws.Range("D13") "19.99"
$ws.Range("E13").Value = "  +0.42%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "7.71"
$ws.Range("E14").Value = "  -1.53%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.362.00"
$ws.Range("E15").Value = "  +3.25%  "

# Row 16
$ws.Range("D16").Value = "2.903.63"
$ws.Range("E16").Value = "  +3.31%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.982"
$ws.Range("E17").Value = "  +1.43%  "

# Row 18
$ws.Range("D18").Value = "51.927.19"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.55"
$ws.Range("E19").Value = "  -1.25%  "

# Row 20
$ws.Range("E20").Value = "  -1.94%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.04"
$ws.Range("E21").Value = "  +3.00%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
Set-TextValue $ws.Range("D23") "70.93"
$ws.Range("E23").Value = "  +0.47%  "

# Row 24
Set-TextValue $ws.Range("D24") "267.94"
$ws.Range("E24").Value = "  -0.70%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.78"
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("E26").Value = "  +8.66%  "

# Row 27
Set-TextValue $ws.Range("D27") "26.76"
$ws.Range("E27").Value = "  +1.85%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
Set-TextValue $ws.Range("D29") "6.92"
$ws.Range("E29").Value = "  +11.95%  "

# Row 30
Set-TextValue $ws.Range("D30") "10.61"
$ws.Range("E30").Value = "  +0.62%  "

# Row 31
$ws.Range("E31").Value = "  +12.18%  "

# Row 32
Set-TextValue $ws.Range("D32") "36.90"
$ws.Range("E32").Value = "  -4.85%  "

# Row 33
Set-TextValue $ws.Range("D33") "5.98"
$ws.Range("E33").Value = "  +5.13%  "

# Row 34
Set-TextValue $ws.Range("D34") "52.88"
$ws.Range("E34").Value = "  +0.21%  "

# Row 35
Set-TextValue $ws.Range("D35") "2.06"
$ws.Range("E35").Value = "  -9.72%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0452"
$ws.Range("E36").Value = "  -0.07%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.998"
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.35"
$ws.Range("E38").Value = "  +4.16%  "

# Row 39
Set-TextValue $ws.Range("D39") "18.56"
$ws.Range("E39").Value = "  -1.78%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.03"
$ws.Range("E40").Value = "  +0.79%  "

# Row 41
Set-TextValue $ws.Range("D41") "2.66"
$ws.Range("E41").Value = "  +3.75%  "

# Row 42
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
Set-TextValue $ws.Range("D43") "22.81"
$ws.Range("E43").Value = "  +3.28%  "

# Row 44
$ws.Range("E44").Value = "  -2.23%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.185.29"
$ws.Range("E45").Value = "  +2.11%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D46") "2.51"
$ws.Range("E46").Value = "  +2.33%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.48"
$ws.Range("E47").Value = "  -1.23%  "

# Row 48
Set-TextValue $ws.Range("D48") "110.99"
$ws.Range("E48").Value = "  -8.66%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.248"
$ws.Range("E49").Value = "  +11.12%  "

# Row 50
$ws.Range("E50").Value = "  +5.59%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.945"
$ws.Range("E51").Value = "  -7.86%  "

